$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Periodo Mora" period value for all workers from 2507 to 2508.
$ws.Range("E16").Value = "2508"
$ws.Range("E17").Value = "2508"
$ws.Range("E18").Value = "2508"

# Update the "Valor Mora" amount for the second worker (row 17).
$ws.Range("G17").Value = 1423500
